# Generate Report for Handback
# Refresh the timestamp values recorded on the handback status report.
#
#   Overview!G2  "Latest HO Xliff Generate Date" for 57812fd9-...md
#   zh-cn!H2     "Correspond Handoff Datetime"   for 57812fd9-...md
#   zh-cn!K2     "Correspond Handback DateTime"  for 57812fd9-...md
#   de-de!K2     "Correspond Handback DateTime"  for 57812fd9-...md

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 05:03:54"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 05:03:49"
$wsZhCn.Range("K2").Value = "2016-08-25 05:04:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-25 05:04:22"
